# Generate Report for Handoff
# - Drops the "75f0a5ad...md / Handoff transform failed" row (row 3) from every
#   sheet, shifting the ".localization-config" row up.
# - Renames the source markdown file (485fc82b... -> 57ba4bf9...) everywhere.
# - Renames the generated xlf handoff files + refreshes their handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldFile = "485fc82b-2dae-42f7-9c4b-467781b55b6e"
$newFile = "57ba4bf9-b79d-4b08-83cd-3fa31246ea3c"
$oldZhXlf = "485fc82b-2dae-42f7-9c4b-467781b55b6e.24e938bce0830702ffd12cabd4404563749daaf7.zh-cn.xlf"
$newZhXlf = "57ba4bf9-b79d-4b08-83cd-3fa31246ea3c.b774cf0b864dccba6fe53cfa7bbadfd6dd8e254a.zh-cn.xlf"
$oldDeXlf = "485fc82b-2dae-42f7-9c4b-467781b55b6e.24e938bce0830702ffd12cabd4404563749daaf7.de-de.xlf"
$newDeXlf = "57ba4bf9-b79d-4b08-83cd-3fa31246ea3c.b774cf0b864dccba6fe53cfa7bbadfd6dd8e254a.de-de.xlf"

$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/07f011b805d225d0b7c3fe5136d46b1984ae0b40/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/07f011b805d225d0b7c3fe5136d46b1984ae0b40/.localization-config"
$zhXlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b5247ea42df662af48a1cbbefca420aa26c6ebed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/"
$deXlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f1ecd1f3698d9ed2d363385be310757060a757d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Hyperlinks.Delete()
$ws1.Rows.Item(3).Delete()

$ws1.Range("A2").Value = "$newFile.md"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$mdUrlBase$newFile.md", [System.Type]::Missing, [System.Type]::Missing, "$newFile.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $cfgUrl, [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

$ws2.Range("A2").Value = "$newFile.md"
$ws2.Range("C2").Value = $newZhXlf
$ws2.Range("D2").Value = "2016-01-28 05:46:20"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$mdUrlBase$newFile.md", [System.Type]::Missing, [System.Type]::Missing, "$newFile.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhXlfUrlBase$newZhXlf", [System.Type]::Missing, [System.Type]::Missing, $newZhXlf) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $cfgUrl, [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Delete()

$ws3.Range("A2").Value = "$newFile.md"
$ws3.Range("C2").Value = $newDeXlf
$ws3.Range("D2").Value = "2016-01-28 05:46:30"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$mdUrlBase$newFile.md", [System.Type]::Missing, [System.Type]::Missing, "$newFile.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deXlfUrlBase$newDeXlf", [System.Type]::Missing, [System.Type]::Missing, $newDeXlf) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $cfgUrl, [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null
